# "added the answer pdf" -- replace the three question cells in column B
# with the machine-learning Q&A content, and widen column B to fit the
# longer text (column B was already best-fit sized to its old content).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "WHAT IS MACHINE LEARNING?"
$ws.Range("B2").Value = "What are the main types of machine learning algorithms ?"
$ws.Range("B3").Value = "Explain the concept of machine learning."

# Column B needs to grow to fit the longest new question text.
$ws.Columns.Item(2).ColumnWidth = 51
